$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 323 (shifts 323:346 down to 324:347).
# Excel copies the formatting of the row above, which already matches what we need.
$ws.Rows(323).Insert()

# Populate the newly inserted row 323 with a duplicate of the (now shifted) row 324
# data, adjusted per the weekly update (date + price range + pieces count).
$ws.Range("A323").Value = 6
$ws.Range("B323").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C323").Value = "Metropolitana"
$ws.Range("D323").Value = 44610
$ws.Range("E323").Value = 13
$ws.Range("F323").Value = 100112032
$ws.Range("G323").Value = "Zapallo italiano"
$ws.Range("H323").Value = "Sin especificar"
$ws.Range("I323").Value = "Primera"
$ws.Range("J323").Value = 400
$ws.Range("K323").Value = 8000
$ws.Range("L323").Value = 9000
$ws.Range("M323").Value = 8575
$ws.Range("N323").Value = "`$/caja 50 unidades"
$ws.Range("O323").Value = "Región Metropolitana"
$ws.Range("P323").Value = 172
$ws.Range("Q323").Value = 50
$ws.Range("R323").Value = "Hortaliza"
